$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force column D (Price) to text format so numeric-looking strings
# like "1.014" or "0.00001040" are preserved exactly as text,
# matching the inlineStr cells in the source workbook.
$ws.Range("D2:D51").NumberFormat = "@"

$ws.Range("D2").Value = '29.552.34'
$ws.Range("E2").Value = '  +0.16%  '
$ws.Range("D3").Value = '1.919.94'
$ws.Range("E3").Value = '  +0.47%  '
$ws.Range("D4").Value = '1.014'
$ws.Range("E4").Value = '  +0.79%  '
$ws.Range("D5").Value = '326.18'
$ws.Range("E5").Value = '  -0.02%  '
$ws.Range("D6").Value = '1.012'
$ws.Range("E6").Value = '  +0.78%  '
$ws.Range("D7").Value = '0.4816'
$ws.Range("E7").Value = '  -0.56%  '
$ws.Range("D8").Value = '0.4055'
$ws.Range("E8").Value = '  -0.55%  '
$ws.Range("D9").Value = '0.08241'
$ws.Range("E9").Value = '  +0.96%  '
$ws.Range("D10").Value = '1.010'
$ws.Range("E10").Value = '  -0.15%  '
$ws.Range("D11").Value = '23.37'
$ws.Range("E11").Value = '  -0.40%  '
$ws.Range("D12").Value = '1.909.27'
$ws.Range("E12").Value = '  -0.98%  '
$ws.Range("D13").Value = '6.064'
$ws.Range("E13").Value = '  +0.67%  '
$ws.Range("D14").Value = '7.249'
$ws.Range("E14").Value = '  +1.92%  '
$ws.Range("D15").Value = '91.70'
$ws.Range("E15").Value = '  +1.34%  '
$ws.Range("D16").Value = '0.06881'
$ws.Range("E16").Value = '  +1.34%  '
$ws.Range("D17").Value = '1.013'
$ws.Range("E17").Value = '  +0.78%  '
$ws.Range("D18").Value = '0.00001040'
$ws.Range("E18").Value = '  -0.08%  '
$ws.Range("D19").Value = '17.57'
$ws.Range("E19").Value = '  -0.87%  '
$ws.Range("D20").Value = '1.011'
$ws.Range("E20").Value = '  +0.73%  '
$ws.Range("D21").Value = '29.565.41'
$ws.Range("E21").Value = '  +0.18%  '
$ws.Range("D22").Value = '5.680'
$ws.Range("E22").Value = '  +1.23%  '
$ws.Range("D23").Value = '11.88'
$ws.Range("E23").Value = '  +0.93%  '
$ws.Range("D24").Value = '2.194'
$ws.Range("E24").Value = '  +1.22%  '
$ws.Range("D25").Value = '2.165.84'
$ws.Range("E25").Value = '  +0.29%  '
$ws.Range("D26").Value = '6.537'
$ws.Range("E26").Value = '  +3.46%  '
$ws.Range("D27").Value = '156.14'
$ws.Range("E27").Value = '  +1.04%  '
$ws.Range("D28").Value = '20.03'
$ws.Range("E28").Value = '  -0.38%  '
$ws.Range("D29").Value = '2.099'
$ws.Range("E29").Value = '  -0.35%  '
$ws.Range("D30").Value = '120.65'
$ws.Range("E30").Value = '  +0.83%  '
$ws.Range("D31").Value = '1.019'
$ws.Range("E31").Value = '  -1.06%  '
$ws.Range("D32").Value = '0.09633'
$ws.Range("E32").Value = '  +0.67%  '
$ws.Range("D33").Value = '5.629'
$ws.Range("E33").Value = '  +1.38%  '
$ws.Range("D34").Value = '3.560'
$ws.Range("E34").Value = '  +0.26%  '
$ws.Range("D35").Value = '1.377'
$ws.Range("E35").Value = '  -1.29%  '
$ws.Range("D36").Value = '0.06346'
$ws.Range("E36").Value = '  +3.85%  '
$ws.Range("D37").Value = '0.02287'
$ws.Range("E37").Value = '  +0.74%  '
$ws.Range("D38").Value = '1.198'
$ws.Range("E38").Value = '  +2.18%  '
$ws.Range("D39").Value = '0.5944'
$ws.Range("E39").Value = '  -0.01%  '
$ws.Range("D40").Value = '10.71'
$ws.Range("E40").Value = '  -0.89%  '
$ws.Range("D41").Value = '7.917'
$ws.Range("E41").Value = '  -0.44%  '
$ws.Range("D42").Value = '0.1848'
$ws.Range("E42").Value = '  -0.55%  '
$ws.Range("D43").Value = '2.466'
$ws.Range("E43").Value = '  +0.06%  '
$ws.Range("D44").Value = '1.282'
$ws.Range("E44").Value = '  +0.05%  '
$ws.Range("D45").Value = '12.44'
$ws.Range("E45").Value = '  -0.24%  '
$ws.Range("D46").Value = '0.07476'
$ws.Range("E46").Value = '  -3.22%  '
$ws.Range("D47").Value = '0.5565'
$ws.Range("E47").Value = '  -0.20%  '
$ws.Range("D48").Value = '1.943'
$ws.Range("E48").Value = '  -0.45%  '
$ws.Range("D49").Value = '118.80'
$ws.Range("E49").Value = '  +3.21%  '
$ws.Range("D50").Value = '2.436'
$ws.Range("E50").Value = '  +3.70%  '
$ws.Range("D51").Value = '72.15'
$ws.Range("E51").Value = '  -0.91%  '
